# This script rebuilds the "Cluster name" / "Active cases" table on Sheet1
# to match the target state described in the commit diff:
#   - several cluster names were removed
#   - several new cluster names were added
#   - the table is re-sorted alphabetically by cluster name
#   - the "Active cases" numbers were updated
#
# Data below is the full set of (ClusterName, ActiveCases) rows in their
# final order, taken directly from the target OOXML in the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,("3035 Campbell Place Aged Care Glen Waverley", 13)
    ,("3364 Assisi Centre Aged Care Rosanna", 20)
    ,("3528 Ottoman Village Aged Care Broadmeadows", 12)
    ,("3622 Olivet Care Aged Care Services Ringwood", 13)
    ,("3633 Lifeview Emerald Glades Aged Care Emerald", 17)
    ,("3652 Regis Aged Care Dandenong North", 12)
    ,("3824 Estia Health South Morang", 17)
    ,("3961 Heritage Care Water Gardens Aged Care Facility Sydenham", 20)
    ,("AG Industries Pty Ltd Thomastown", 13)
    ,("Aintree Primary School Aintree", 18)
    ,("Australian Meat Group Abattoir Dandenong South", 20)
    ,("Bacchus Marsh Childcare and Kindergarten Centre Bacchus Marsh", 22)
    ,("Baden Powell College Tarneit", 12)
    ,("Bandiana Primary School Bandiana", 10)
    ,("CREST Children's Sanctuary Dandenong", 10)
    ,("Dandenong South Primary School Dandenong", 10)
    ,("Hamlyn Views School Hamlyn Heights", 12)
    ,("Hazelwood North Primary School Hazelwood", 10)
    ,("Hippity Hop Childcare and Kindergarten Pakenham", 10)
    ,("Lilydale Motor Inn Lilydale", 10)
    ,("Lowanna College Newborough", 12)
    ,("M.C. Herd Corio", 10)
    ,("Master Poultry Group West Footscray", 10)
    ,("Morwell Park Primary School Morwell", 17)
    ,("Northern Bay College Goldsworthy 9-12 Campus Corio", 16)
    ,("Northern Bay College Wexford Campus Corio", 54)
    ,("Oakleigh South Primary School Oakleigh South", 11)
    ,("Rosewood Downs Special Accommodation Home Dandenong", 14)
    ,("Saint Augustines Primary School Wodonga", 10)
    ,("Saint Monica's Primary School Wodonga", 11)
    ,("Smartie Pants Early Learning and Development Diamond Creek", 10)
    ,("St Josephs Catholic Primary School Warragul", 12)
    ,("St Mary's Primary School Swan Hill", 10)
    ,("St. Brendans Catholic Primary School Lakes Entrance", 10)
    ,("Stockdale Road Primary School Traralgon", 15)
    ,("TUROSI PTY LTD Thomastown", 10)
    ,("Wodonga Cemetery Wodonga", 35)
    ,("Wodonga Primary School Wodonga", 24)
    ,("Wodonga Senior Secondary College Wodonga", 24)
    ,("Wodonga South Primary School Wodonga", 39)
    ,("Woodend Primary School Woodend", 20)
    ,("Yallourn Power Station Yallourn", 10)
    ,("Yooralla Disability Residential Care St Albans", 11)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
